$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 140, pushing existing rows 140-142 down to 141-143
$ws.Rows("140:140").Insert()

# Copy the date cell style (column D) from the row below (now row 141, previously row 140)
# so the new row's date cell matches the existing date formatting
$ws.Range("D141").Copy()
$ws.Range("D140").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Populate the new row 140 with the new record's data
$ws.Range("A140").Value = 10
$ws.Range("B140").Value = "Vega Modelo de Temuco"
$ws.Range("C140").Value = "La Araucanía"
$ws.Range("D140").Value = 45239
$ws.Range("E140").Value = 9
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100101
$ws.Range("H140").Value = "Berries"
$ws.Range("I140").Value = 100101001
$ws.Range("J140").Value = "Arándano (blue)"
$ws.Range("K140").Value = "Sin especificar"
$ws.Range("L140").Value = "Primera"
$ws.Range("M140").Value = 250
$ws.Range("N140").Value = 5500
$ws.Range("O140").Value = 5500
$ws.Range("P140").Value = 5500
$ws.Range("Q140").Value = "$/kilo"
$ws.Range("R140").Value = "Región del Maule"
$ws.Range("S140").Value = 5500
$ws.Range("T140").Value = 1
